$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1517.6316
$ws.Range("I28").Value = 967.8570999999999
$ws.Range("J28").Value = 3057
$ws.Range("K28").Value = 967.8570999999999
$ws.Range("L28").Value = 3057
$ws.Range("M28").Value = -482.8570999999999
$ws.Range("N28").Value = -4027
$ws.Range("H32").Value = 1034.4
$ws.Range("I32").Value = 990
$ws.Range("J32").Value = 1078.8
$ws.Range("K32").Value = 990
$ws.Range("L32").Value = 1078.8
$ws.Range("M32").Value = -664
$ws.Range("N32").Value = -1730.8
$ws.Range("H33").Value = 28571894
$ws.Range("I33").Value = 192.90909
$ws.Range("J33").Value = 76924000
$ws.Range("K33").Value = 192.90909
$ws.Range("L33").Value = 76924000
$ws.Range("M33").Value = 36.09091000000001
$ws.Range("N33").Value = -76924458
$ws.Range("H41").Value = 351.6111
$ws.Range("I41").Value = 383
$ws.Range("J41").Value = 270
$ws.Range("K41").Value = 383
$ws.Range("L41").Value = 270
$ws.Range("M41").Value = 57
$ws.Range("N41").Value = -1150
$ws.Range("H62").Value = 4193
$ws.Range("I62").Value = 3811
$ws.Range("J62").Value = 4359.087
$ws.Range("K62").Value = 3811
$ws.Range("L62").Value = 4359.087
$ws.Range("M62").Value = -3187
$ws.Range("N62").Value = -5607.087
$ws.Range("H65").Value = 4193
$ws.Range("I65").Value = 3811
$ws.Range("J65").Value = 4359.087
$ws.Range("K65").Value = 19055
$ws.Range("L65").Value = 21795.435
$ws.Range("M65").Value = -15935
$ws.Range("N65").Value = -28035.435
$ws.Range("H76").Value = 3008.3333
$ws.Range("I76").Value = 3008.3333
$ws.Range("K76").Value = 3008.3333
$ws.Range("M76").Value = -2693.3333
$ws.Range("H79").Value = 3008.3333
$ws.Range("I79").Value = 3008.3333
$ws.Range("K79").Value = 3008.3333
$ws.Range("M79").Value = -1916.3333
$ws.Range("H92").Value = 955.64703
$ws.Range("I92").Value = 718.9231
$ws.Range("K92").Value = 718.9231
$ws.Range("M92").Value = 529.0769
$ws.Range("H98").Value = 1780.72
$ws.Range("I98").Value = 1286.1666
$ws.Range("J98").Value = 2237.2307
$ws.Range("K98").Value = 1286.1666
$ws.Range("L98").Value = 2237.2307
$ws.Range("M98").Value = 211.8334
$ws.Range("N98").Value = -5233.2307
$ws.Range("H100").Value = 2232.353
$ws.Range("I100").Value = 1306.25
$ws.Range("J100").Value = 3055.5557
$ws.Range("K100").Value = 1306.25
$ws.Range("L100").Value = 3055.5557
$ws.Range("M100").Value = -765.25
$ws.Range("N100").Value = -4137.5557
$ws.Range("H122").Value = 1780.72
$ws.Range("I122").Value = 1286.1666
$ws.Range("J122").Value = 2237.2307
$ws.Range("K122").Value = 3858.4998
$ws.Range("L122").Value = 6711.6921
$ws.Range("M122").Value = -1408.4998
$ws.Range("N122").Value = -11611.6921
$ws.Range("H135").Value = 21277046
$ws.Range("I135").Value = 426.06668
$ws.Range("J135").Value = 500001000
$ws.Range("K135").Value = 3834.60012
$ws.Range("L135").Value = 4500009000
$ws.Range("M135").Value = -1299.60012
$ws.Range("N135").Value = -4500014070
$ws.Range("H138").Value = 2064392.4
$ws.Range("I138").Value = 1005.6964
$ws.Range("J138").Value = 4882676.5
$ws.Range("K138").Value = 3017.0892
$ws.Range("L138").Value = 14648029.5
$ws.Range("M138").Value = 2122.9108
$ws.Range("N138").Value = -14658309.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4978.9297
$ws.Range("I32").Value = 3921.8936
$ws.Range("J32").Value = 9947
$ws.Range("K32").Value = 3921.8936
$ws.Range("L32").Value = 9947
$ws.Range("M32").Value = -3634.8936
$ws.Range("N32").Value = -10521
$ws.Range("H52").Value = 17999.5
$ws.Range("J52").Value = 17999.5
$ws.Range("L52").Value = 17999.5
$ws.Range("N52").Value = -18635.5
$ws.Range("H110").Value = 1277.5385
$ws.Range("I110").Value = 1354.0952
$ws.Range("J110").Value = 956
$ws.Range("K110").Value = 1354.0952
$ws.Range("L110").Value = 956
$ws.Range("M110").Value = 690.9048
$ws.Range("N110").Value = -5046

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 29900
$ws.Range("J51").Value = 29900
$ws.Range("L51").Value = 29900
$ws.Range("N51").Value = -30882
$ws.Range("H57").Value = 37200
$ws.Range("J57").Value = 37200
$ws.Range("L57").Value = 37200
$ws.Range("N57").Value = -38640
$ws.Range("H107").Value = 1463.3438
$ws.Range("I107").Value = 1365.2916
$ws.Range("J107").Value = 1757.5
$ws.Range("K107").Value = 1365.2916
$ws.Range("L107").Value = 1757.5
$ws.Range("M107").Value = 554.7084
$ws.Range("N107").Value = -5597.5
$ws.Range("H136").Value = 37200
$ws.Range("J136").Value = 37200
$ws.Range("L136").Value = 37200
$ws.Range("N136").Value = -47400

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 792.75
$ws.Range("I16").Value = 640.1667
$ws.Range("J16").Value = 945.3333
$ws.Range("K16").Value = 640.1667
$ws.Range("L16").Value = 945.3333
$ws.Range("M16").Value = -353.1667
$ws.Range("N16").Value = -1519.3333
$ws.Range("H58").Value = 5177.875
$ws.Range("I58").Value = 5337.7827
$ws.Range("J58").Value = 1500
$ws.Range("K58").Value = 5337.7827
$ws.Range("L58").Value = 1500
$ws.Range("M58").Value = -5134.7827
$ws.Range("N58").Value = -1906
$ws.Range("H113").Value = 792.75
$ws.Range("I113").Value = 640.1667
$ws.Range("J113").Value = 945.3333
$ws.Range("K113").Value = 640.1667
$ws.Range("L113").Value = 945.3333
$ws.Range("M113").Value = 1529.8333
$ws.Range("N113").Value = -5285.3333
$ws.Range("H132").Value = 824769.1
$ws.Range("I132").Value = 1701.8918
$ws.Range("J132").Value = 4631455
$ws.Range("K132").Value = 5105.6754
$ws.Range("L132").Value = 13894365
$ws.Range("M132").Value = -2575.6754
$ws.Range("N132").Value = -13899425
$ws.Range("H134").Value = 33335358
$ws.Range("I134").Value = 4168754.5
$ws.Range("J134").Value = 111112970
$ws.Range("K134").Value = 12506263.5
$ws.Range("L134").Value = 333338910
$ws.Range("M134").Value = -12503728.5
$ws.Range("N134").Value = -333343980
$ws.Range("H136").Value = 5177.875
$ws.Range("I136").Value = 5337.7827
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 16013.3481
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -13463.3481
$ws.Range("N136").Value = -9600

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 2409.2222
$ws.Range("J115").Value = 3100
$ws.Range("L115").Value = 9300
$ws.Range("N115").Value = -11650

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1279.9
$ws.Range("I113").Value = 1185.5714
$ws.Range("K113").Value = 1185.5714
$ws.Range("M113").Value = 984.4286
$ws.Range("H122").Value = 51960.523
$ws.Range("I122").Value = 86429.164
$ws.Range("J122").Value = 6002.3335
$ws.Range("K122").Value = 259287.492
$ws.Range("L122").Value = 18007.0005
$ws.Range("M122").Value = -256837.492
$ws.Range("N122").Value = -22907.0005

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1710.8
$ws.Range("I61").Value = 1318
$ws.Range("K61").Value = 1318
$ws.Range("M61").Value = -1116
$ws.Range("H113").Value = 1710.8
$ws.Range("I113").Value = 1318
$ws.Range("K113").Value = 1318
$ws.Range("M113").Value = 852

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 442.45456
$ws.Range("I113").Value = 302.88
$ws.Range("J113").Value = 878.625
$ws.Range("K113").Value = 908.64
$ws.Range("L113").Value = 2635.875
$ws.Range("M113").Value = 1261.36
$ws.Range("N113").Value = -6975.875
$ws.Range("H132").Value = 3075
$ws.Range("I132").Value = 3101.8572
$ws.Range("J132").Value = 2999.8
$ws.Range("K132").Value = 9305.571599999999
$ws.Range("L132").Value = 8999.400000000001
$ws.Range("M132").Value = -6775.571599999999
$ws.Range("N132").Value = -14059.4
